$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the second "Periodo Mora" row (period 2505) for the worker,
# shifting the signature block up (rows 22/23 -> 21/22).
$ws.Rows("17:17").Delete()

# Update the summary totals to reflect the remaining single period:
# "VALOR MORA" total is now just the remaining row's Valor Mora (36062
# instead of 36062 + 7592 = 43654).
$ws.Range("E11").Value = 36062

# "Cant. Periodos" drops from 2 to 1.
$ws.Range("F13").Value = 1
